$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New column header "PT%" in C1 (adds a new shared string) ---
$ws.Range("C1").Value = "PT%"

# --- 2. New "PT%" data values for rows 3-16 (column C) ---
$ws.Range("C3").Value = 73.4
$ws.Range("C4").Value = 65.5
$ws.Range("C5").Value = 67.7
$ws.Range("C6").Value = 61.4
$ws.Range("C7").Value = 76.5
$ws.Range("C8").Value = 88.4
$ws.Range("C9").Value = 93
$ws.Range("C10").Value = 98.7
$ws.Range("C11").Value = 95.3
$ws.Range("C12").Value = 98.7
$ws.Range("C13").Value = 93.7
$ws.Range("C14").Value = 89
$ws.Range("C15").Value = 92.1
$ws.Range("C16").Value = 101.2

# --- 3. Column C width widened to fit the new header/data (stored width 17) ---
$ws.Columns.Item(3).ColumnWidth = 16.285714285714285

# --- 4. New "[b,value]" helper formula for D20 (row 20 mirrors row 2) ---
$ws.Range("D20").Formula = '="["&$B2&","&D2&"]"'

# --- 5. New "[b,value]" helper formulas in columns C and D for rows 21-34
#         (each row N mirrors data row N-18), plus column D only for rows 35-36
#         (rows 17-18 have no PT% data in column C) ---
for ($r = 21; $r -le 34; $r++) {
    $src = $r - 18
    $ws.Range("C$r").Formula = '="["&$B' + $src + '&","&C' + $src + '&"]"'
    $ws.Range("D$r").Formula = '="["&$B' + $src + '&","&D' + $src + '&"]"'
}
for ($r = 35; $r -le 36; $r++) {
    $src = $r - 18
    $ws.Range("D$r").Formula = '="["&$B' + $src + '&","&D' + $src + '&"]"'
}

# --- 6. Update the sheet view's scroll position and selection ---
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("F21:F34").Select()
